$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("I9").Formula = "=OFFSET(A1, -1, 0)"
$ws.Range("G10").Formula = "=INDEX(A1:E5, 6, 1)"
$ws.Range("I10").Formula = "=OFFSET(A1, 0, -1)"
$ws.Range("M10").Formula = "=INDIRECT(K5)"
$ws.Range("G11").Formula = "=INDEX(A1:E5, 1, 6)"
$ws.Range("M11").Formula = "=INDIRECT("""")"
